# Correct the loan-product code in both sheets' "productname" value cell
# (B1): insert the missing dash right after "246" so it reads
# "246-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME" instead of
# "246MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME".

$wb = $excel.ActiveWorkbook

$correctedName = "246-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsInput.Range("B1").Value = $correctedName
$wsInput.Range("B1").Select()

$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")
$wsOutput.Range("B1").Value = $correctedName

# Make the output sheet the active tab/selection, leaving the input sheet's
# selection on B1 (matching the saved view state of the workbook).
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
